$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression
$ws.Range("B2").Value = 1802482372333881
$ws.Range("C2").Value = 1802482372333881
$ws.Range("D2").Value = 1802482372333879

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.01633735883174511
$ws.Range("C3").Value = 0.01661492937183706
$ws.Range("D3").Value = 2699500889566.902

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01965645267219591
$ws.Range("C4").Value = 0.02026980981117747
$ws.Range("D4").Value = 0.05948685421572617

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 9812119808140.598
$ws.Range("C5").Value = 9822558771007.205
$ws.Range("D5").Value = 78938700645022.83
